# Update automàtic: dades i banners [2026-02-24 18:50]
# Applies the refreshed MeteoCat extraction timestamps and sensor readings
# recorded for each station row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = '2026-02-24 18:48:23'
$ws.Range("E3").Value = '2026-02-24 18:48:25'
$ws.Range("L3").Value = '26.6 km/h - 126º 18:25 TU'
$ws.Range("E4").Value = '2026-02-24 18:48:27'
$ws.Range("J4").Value = '1019.9 hPa'
$ws.Range("E5").Value = '2026-02-24 18:48:30'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '27%'
$ws.Range("O5").Value = '6.1 °C'
$ws.Range("E6").Value = '2026-02-24 18:48:32'
$ws.Range("J6").Value = '1019.8 hPa'
$ws.Range("E7").Value = '2026-02-24 18:48:35'
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '69%'
$ws.Range("O7").Value = '14.3 °C'
$ws.Range("E8").Value = '2026-02-24 18:48:37'
$ws.Range("E9").Value = '2026-02-24 18:48:40'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '79%'
$ws.Range("E10").Value = '2026-02-24 18:48:42'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '75%'
$ws.Range("O10").Value = '11.7 °C'
$ws.Range("E11").Value = '2026-02-24 18:48:45'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '68%'
$ws.Range("E12").Value = '2026-02-24 18:48:47'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '90%'
$ws.Range("E13").Value = '2026-02-24 18:48:49'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '62%'
$ws.Range("J13").Value = '1023.5 hPa'
$ws.Range("L13").Value = '18.4 km/h - 122º 18:29 TU'
$ws.Range("O13").Value = '6.7 °C'
$ws.Range("E14").Value = '2026-02-24 18:48:52'
$ws.Range("E15").Value = '2026-02-24 18:48:54'
$ws.Range("E16").Value = '2026-02-24 18:48:56'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '19%'
$ws.Range("N16").Value = '2.1 °C 18:29 TU'
$ws.Range("E17").Value = '2026-02-24 18:48:59'
$ws.Range("E18").Value = '2026-02-24 18:49:01'
$ws.Range("J18").Value = '1020.3 hPa'
$ws.Range("O18").Value = '11.3 °C'
$ws.Range("E19").Value = '2026-02-24 18:49:04'
$ws.Range("E20").Value = '2026-02-24 18:49:06'
$ws.Range("O20").Value = '3.4 °C'
$ws.Range("E21").Value = '2026-02-24 18:49:09'
$ws.Range("J21").Value = '1022.2 hPa'
$ws.Range("O21").Value = '9.7 °C'
$ws.Range("E22").Value = '2026-02-24 18:49:11'
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = '22%'
$ws.Range("O22").Value = '3.7 °C'
$ws.Range("E23").Value = '2026-02-24 18:49:13'
$ws.Range("E24").Value = '2026-02-24 18:49:16'
$ws.Range("J24").Value = '1021.5 hPa'
$ws.Range("E25").Value = '2026-02-24 18:49:18'
$ws.Range("N25").Value = '3.8 °C 18:29 TU'
$ws.Range("O25").Value = '7.1 °C'
$ws.Range("E26").Value = '2026-02-24 18:49:21'
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = '41%'
$ws.Range("J26").Value = '1018.6 hPa'
$ws.Range("E27").Value = '2026-02-24 18:49:23'
$ws.Range("E28").Value = '2026-02-24 18:49:26'
$ws.Range("J28").Value = '1020.2 hPa'
$ws.Range("E29").Value = '2026-02-24 18:49:28'
$ws.Range("E30").Value = '2026-02-24 18:49:31'
$ws.Range("J30").Value = '1019.9 hPa'
$ws.Range("O30").Value = '13.4 °C'
$ws.Range("E31").Value = '2026-02-24 18:49:33'
$ws.Range("E32").Value = '2026-02-24 18:49:36'
$ws.Range("E33").Value = '2026-02-24 18:49:38'
$ws.Range("J33").Value = '1021.8 hPa'
$ws.Range("O33").Value = '8.6 °C'
$ws.Range("E34").Value = '2026-02-24 18:49:41'
$ws.Range("E35").Value = '2026-02-24 18:49:43'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '37%'
$ws.Range("J35").Value = '1020.6 hPa'
$ws.Range("E36").Value = '2026-02-24 18:49:46'
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '78%'
$ws.Range("J36").Value = '1020.1 hPa'
$ws.Range("E37").Value = '2026-02-24 18:49:48'
$ws.Range("E38").Value = '2026-02-24 18:49:50'
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = '70%'
$ws.Range("E39").Value = '2026-02-24 18:49:53'
$ws.Range("O39").Value = '4.7 °C'
$ws.Range("E40").Value = '2026-02-24 18:49:55'
$ws.Range("J40").Value = '1022.8 hPa'
$ws.Range("E41").Value = '2026-02-24 18:49:58'
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = '77%'
$ws.Range("J41").Value = '1020.7 hPa'
$ws.Range("E42").Value = '2026-02-24 18:50:00'
$ws.Range("E43").Value = '2026-02-24 18:50:02'
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = '70%'
$ws.Range("O43").Value = '10.7 °C'
$ws.Range("E44").Value = '2026-02-24 18:50:05'
$ws.Range("E45").Value = '2026-02-24 18:50:07'
$ws.Range("J45").Value = '1020.9 hPa'
$ws.Range("L45").Value = '21.2 km/h - 111º 18:02 TU'
$ws.Range("E46").Value = '2026-02-24 18:50:10'
$ws.Range("J46").Value = '1021.5 hPa'
$ws.Range("O46").Value = '10.5 °C'
